# Add a new closing "Questions?" slide at the end of the deck.
#
# The new slide reuses the Title-Slide layout (ctrTitle / subTitle
# placeholders) already used by slide 1 ("Design Tips"), so the cleanest
# and most faithful way to reproduce it is to duplicate slide 1, move the
# duplicate to the end of the deck, and replace its title text.

$p = $ppt.ActivePresentation

# Duplicate the first slide (title-slide layout: ctrTitle + subTitle).
$firstSlide = $p.Slides.Item(1)
$dup = $firstSlide.Duplicate()
$newSlide = $dup.Item(1)

# Move the duplicate to be the last slide in the deck.
$newSlide.MoveTo($p.Slides.Count)

# Update the title text; leave the subtitle placeholder empty.
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Questions?"
